# Apply the "upload new version with timestamp" update to the low-stock
# report: insert two new low-stock products (CONGESTAL SYRUP 120 ML and
# GAST-REG 200 MG 30 TABS.) into the sorted product table, renumber the
# row index column, refresh the grand-total cell, and bump the generated
# timestamp shown in the footer.

$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

# --- 1. Grow the table by two rows -----------------------------------
# The table currently spans rows 7-17 (11 products) followed by the total
# row (18) and the footer row (19). Insert two blank rows right before the
# total row so the table becomes rows 7-19 (13 products), pushing the total
# row to 20 and the footer row to 21.
$ws.Rows.Item(18).Insert()
$ws.Rows.Item(19).Insert()

# Clone formatting (styles, borders, number formats) from the last existing
# product row (17) into the two new rows.
$ws.Range("A17:Q17").Copy()
$ws.Range("A18:Q18").PasteSpecial(-4122)
$ws.Range("A17:Q17").Copy()
$ws.Range("A19:Q19").PasteSpecial(-4122)
$excel.CutCopyMode = 0

# Re-create the per-row merged cell groups for the two new rows, matching
# the pattern used by every other product row.
$ws.Range("A18:B18").Merge()
$ws.Range("C18:G18").Merge()
$ws.Range("H18:K18").Merge()
$ws.Range("L18:M18").Merge()
$ws.Range("N18:O18").Merge()
$ws.Range("A19:B19").Merge()
$ws.Range("C19:G19").Merge()
$ws.Range("H19:K19").Merge()
$ws.Range("L19:M19").Merge()
$ws.Range("N19:O19").Merge()

# Restore the alternating row-height rhythm used throughout the table.
$ws.Rows.Item(18).RowHeight = 24.75
$ws.Rows.Item(19).RowHeight = 25.5

# --- 2. Rewrite every product row with the final, sorted data --------
# Columns: A = row number, C = item name, H = current balance,
# L = order limit, N = price, P = sale price, Q = transaction count.
$products = @(
  @{ A = 1;  C = "BABETONE SYRUP SUGAR FREE 120 ML"; H = "0:0";  L = "1"; N = "35.00";  P = "35.0000"; Q = "1:0" },
  @{ A = 2;  C = "CONGESTAL SYRUP 120 ML";            H = "2:0";  L = "1"; N = "44.00";  P = "44.0000"; Q = "1:0" },
  @{ A = 3;  C = "DEPO-PEN 1.2 MIU VIAL.";             H = "3:0";  L = "1"; N = "25.00";  P = "25.0000"; Q = "1:0" },
  @{ A = 4;  C = "ETHOXA 250MG/5ML SYRUP 120ML";       H = "0:0";  L = "1"; N = "99.00";  P = "99.0000"; Q = "1:0" },
  @{ A = 5;  C = "GAST-REG 200 MG 30 TABS.";           H = "0:3";  L = "1"; N = "84.00";  P = "27.7200"; Q = "0:1" },
  @{ A = 6;  C = "KETOLAC 10MG 20 TAB";                H = "0:0";  L = "1"; N = "38.00";  P = "38.0000"; Q = "1:0" },
  @{ A = 7;  C = "TORSERETIC 100MG 30 TABS.";          H = "1:0";  L = "1"; N = "261.00"; P = "86.1300"; Q = "0:1" },
  @{ A = 8;  C = "VOLTAREN 50MG 20 TAB.";              H = "0:1";  L = "1"; N = "48.00";  P = "24.0000"; Q = "0:1" },
  @{ A = 9;  C = "الويز كبير بالاجنحه";                H = "16:0"; L = "0"; N = "60.00";  P = "60.0000"; Q = "1:0" },
  @{ A = 10; C = "جل رويال";                           H = "0:0";  L = "0"; N = "40.00";  P = "40.0000"; Q = "1:0" },
  @{ A = 11; C = "حلق";                                H = "22:0"; L = "0"; N = "10.00";  P = "10.0000"; Q = "1:0" },
  @{ A = 12; C = "شفاط ثدي الجو";                      H = "1:0";  L = "0"; N = "25.00";  P = "25.0000"; Q = "1:0" },
  @{ A = 13; C = "كالونا ";                            H = "0:0";  L = "0"; N = "15.00";  P = "15.0000"; Q = "1:0" }
)

$row = 7
foreach ($p in $products) {
  $ws.Range("A$row").Value = $p.A
  $ws.Range("C$row").Value = $p.C
  $ws.Range("H$row").Value = $p.H
  $ws.Range("L$row").Value = $p.L
  $ws.Range("N$row").Value = $p.N
  $ws.Range("P$row").Value = $p.P
  $ws.Range("Q$row").Value = $p.Q
  $row = $row + 1
}

# --- 3. Refresh the grand total and footer ----------------------------
$ws.Range("P20").Value = 528.85

$ws.Range("A21").Value = "Wednesday, 20 August, 2025 10:59 AM"

Write-Host "Low-stock report updated: 13 products, total 528.85"
